$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 384
$ws1.Range("F4").Value = 2161
$ws1.Range("F5").Value = 63
$ws1.Range("F6").Value = 12539
$ws1.Range("F7").Value = 12539
$ws1.Range("F11").Value = 454
$ws1.Range("F12").Value = 1144
$ws1.Range("F13").Value = 939
$ws1.Range("F14").Value = 13644
$ws1.Range("F15").Value = 13948
$ws1.Range("F23").Value = 474
$ws1.Range("F24").Value = 5040

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 385
$ws4.Range("F4").Value = 2161
$ws4.Range("F5").Value = 63
$ws4.Range("F6").Value = 12539
$ws4.Range("F7").Value = 12539
$ws4.Range("F11").Value = 454
$ws4.Range("F12").Value = 1144
$ws4.Range("F13").Value = 939
$ws4.Range("F14").Value = 13644
$ws4.Range("F15").Value = 13948
$ws4.Range("F23").Value = 474
$ws4.Range("F24").Value = 5040
